$p = $ppt.ActivePresentation

# Remove all slides from the presentation (the deck's sldIdLst becomes empty).
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $p.Slides.Item($i).Delete()
}
